# Apply updated crypto price/volume(1h) data (GitHub Actions refresh).
# Values are written as text (apostrophe-prefixed) to match the workbook's
# original inline-string cell type and avoid Excel auto-converting
# numeric-looking strings (e.g. "514.70") into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '57.146.62'
$ws.Range('E2').Value = "'" + '  -0.48%  '
$ws.Range('D3').Value = "'" + '3.053.18'
$ws.Range('E3').Value = "'" + '  +1.24%  '
$ws.Range('E4').Value = "'" + '  +0.17%  '
$ws.Range('D5').Value = "'" + '514.70'
$ws.Range('E5').Value = "'" + '  +1.19%  '
$ws.Range('D6').Value = "'" + '140.94'
$ws.Range('E6').Value = "'" + '  +0.74%  '
$ws.Range('D7').Value = "'" + '1.00'
$ws.Range('E7').Value = "'" + '  +0.14%  '
$ws.Range('E8').Value = "'" + '  +1.25%  '
$ws.Range('E9').Value = "'" + '  -4.51%  '
$ws.Range('E10').Value = "'" + '  -0.74%  '
$ws.Range('D11').Value = "'" + '0.377'
$ws.Range('E11').Value = "'" + '  +3.12%  '
$ws.Range('D12').Value = "'" + '3.578.48'
$ws.Range('E12').Value = "'" + '  +1.40%  '
$ws.Range('E13').Value = "'" + '  -2.96%  '
$ws.Range('D14').Value = "'" + '27.00'
$ws.Range('E14').Value = "'" + '  +2.44%  '
$ws.Range('E15').Value = "'" + '  +1.62%  '
$ws.Range('D16').Value = "'" + '57.186.62'
$ws.Range('E16').Value = "'" + '  -0.37%  '
$ws.Range('D17').Value = "'" + '6.17'
$ws.Range('E17').Value = "'" + '  -0.85%  '
$ws.Range('D18').Value = "'" + '3.054.88'
$ws.Range('E18').Value = "'" + '  +1.48%  '
$ws.Range('D19').Value = "'" + '13.47'
$ws.Range('E19').Value = "'" + '  +4.99%  '
$ws.Range('D20').Value = "'" + '8.14'
$ws.Range('E20').Value = "'" + '  +2.15%  '
$ws.Range('D21').Value = "'" + '331.48'
$ws.Range('E21').Value = "'" + '  +1.17%  '
$ws.Range('E22').Value = "'" + '  +0.14%  '
$ws.Range('D23').Value = "'" + '0.507'
$ws.Range('E23').Value = "'" + '  +1.71%  '
$ws.Range('D24').Value = "'" + '65.77'
$ws.Range('E24').Value = "'" + '  +1.88%  '
$ws.Range('D25').Value = "'" + '3.176.57'
$ws.Range('E25').Value = "'" + '  +1.29%  '
$ws.Range('E26').Value = "'" + '  +0.15%  '
$ws.Range('E27').Value = "'" + '  -1.10%  '
$ws.Range('D28').Value = "'" + '0.0₃0893'
$ws.Range('E28').Value = "'" + '  -2.93%  '
$ws.Range('D29').Value = "'" + '6.76'
$ws.Range('E29').Value = "'" + '  -0.39%  '
$ws.Range('D30').Value = "'" + '7.18'
$ws.Range('E30').Value = "'" + '  -2.27%  '
$ws.Range('E31').Value = "'" + '  +0.00%  '
$ws.Range('E32').Value = "'" + '  +1.17%  '
$ws.Range('E33').Value = "'" + '  +1.08%  '
$ws.Range('D34').Value = "'" + '4.71'
$ws.Range('E34').Value = "'" + '  -1.36%  '
$ws.Range('D35').Value = "'" + '151.33'
$ws.Range('E35').Value = "'" + '  -1.66%  '
$ws.Range('D36').Value = "'" + '5.95'
$ws.Range('E36').Value = "'" + '  +1.04%  '
$ws.Range('E37').Value = "'" + '  -0.07%  '
$ws.Range('D38').Value = "'" + '25.31'
$ws.Range('E38').Value = "'" + '  +2.14%  '
$ws.Range('D39').Value = "'" + '0.0677'
$ws.Range('E39').Value = "'" + '  +0.00%  '
$ws.Range('B40').Value = "'" + 'OKB'
$ws.Range('C40').Value = "'" + 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = "'" + '36.91'
$ws.Range('E40').Value = "'" + '  -2.69%  '
$ws.Range('B41').Value = "'" + 'Filecoin'
$ws.Range('C41').Value = "'" + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = "'" + '3.90'
$ws.Range('E41').Value = "'" + '  +0.80%  '
$ws.Range('E42').Value = "'" + '  +0.19%  '
$ws.Range('D43').Value = "'" + '0.664'
$ws.Range('E43').Value = "'" + '  +2.11%  '
$ws.Range('E44').Value = "'" + '  -1.20%  '
$ws.Range('D45').Value = "'" + '2.205.00'
$ws.Range('E45').Value = "'" + '  -0.88%  '
$ws.Range('D46').Value = "'" + '6.08'
$ws.Range('E46').Value = "'" + '  +0.34%  '
$ws.Range('D47').Value = "'" + '0.959'
$ws.Range('E47').Value = "'" + '  -2.30%  '
$ws.Range('D48').Value = "'" + '20.34'
$ws.Range('E48').Value = "'" + '  +3.95%  '
$ws.Range('E49').Value = "'" + '  +1.22%  '
$ws.Range('E50').Value = "'" + '  +0.02%  '
$ws.Range('D51').Value = "'" + '0.0171'
$ws.Range('E51').Value = "'" + '  +6.46%  '
